$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds text-formatted numeric-looking strings (e.g. "54.518.88").
# Force the cells to remain Text-formatted so Excel does not reinterpret them as numbers/dates.
$ws.Range("D2:D51").NumberFormat = "@"

# Update cryptocurrency prices (column D) and hourly volume percentages (column E)
$ws.Range("D2").Value = "54.518.88"
$ws.Range("E2").Value = "  +5.08%  "
$ws.Range("D3").Value = "2.419.48"
$ws.Range("E3").Value = "  +4.61%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "471.07"
$ws.Range("E5").Value = "  +8.54%  "
$ws.Range("D6").Value = "136.71"
$ws.Range("E6").Value = "  +12.39%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("D8").Value = "0.497"
$ws.Range("E8").Value = "  +7.25%  "
$ws.Range("D9").Value = "2.425.82"
$ws.Range("E9").Value = "  +4.43%  "
$ws.Range("D10").Value = "5.48"
$ws.Range("E10").Value = "  +7.26%  "
$ws.Range("D11").Value = "0.0944"
$ws.Range("E11").Value = "  +5.84%  "
$ws.Range("E12").Value = "  +5.17%  "
$ws.Range("D13").Value = "0.122"
$ws.Range("E13").Value = "  +2.13%  "
$ws.Range("D14").Value = "2.830.35"
$ws.Range("E14").Value = "  +6.05%  "
$ws.Range("D15").Value = "54.693.30"
$ws.Range("E15").Value = "  +5.54%  "
$ws.Range("D16").Value = "20.10"
$ws.Range("E16").Value = "  +6.77%  "
$ws.Range("D17").Value = "0.0000131"
$ws.Range("E17").Value = "  +10.68%  "
$ws.Range("D18").Value = "2.430.47"
$ws.Range("E18").Value = "  +5.96%  "
$ws.Range("D19").Value = "4.28"
$ws.Range("E19").Value = "  +5.71%  "
$ws.Range("D20").Value = "9.75"
$ws.Range("E20").Value = "  +10.02%  "
$ws.Range("D21").Value = "308.54"
$ws.Range("E21").Value = "  +4.22%  "
$ws.Range("D22").Value = "0.993"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("D23").Value = "5.60"
$ws.Range("E23").Value = "  +8.14%  "
$ws.Range("D24").Value = "56.54"
$ws.Range("E24").Value = "  +6.17%  "
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").Value = "0.396"
$ws.Range("E26").Value = "  +6.72%  "
$ws.Range("E27").Value = "  +18.97%  "
$ws.Range("D28").Value = "2.533.30"
$ws.Range("E28").Value = "  +8.75%  "
$ws.Range("D29").Value = "7.19"
$ws.Range("E29").Value = "  +5.81%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.39%  "
$ws.Range("D31").Value = "0.0₃0752"
$ws.Range("E31").Value = "  +13.05%  "
$ws.Range("D32").Value = "147.94"
$ws.Range("E32").Value = "  +2.71%  "
$ws.Range("D33").Value = "17.93"
$ws.Range("E33").Value = "  +5.85%  "
$ws.Range("E34").Value = "  +9.90%  "
$ws.Range("D35").Value = "5.05"
$ws.Range("E35").Value = "  +5.38%  "
$ws.Range("D36").Value = "1.10"
$ws.Range("E36").Value = "  +11.24%  "
$ws.Range("D37").Value = "3.50"
$ws.Range("E37").Value = "  +5.71%  "
$ws.Range("D38").Value = "0.824"
$ws.Range("E38").Value = "  +8.16%  "
$ws.Range("D39").Value = "33.46"
$ws.Range("E39").Value = "  +4.83%  "
$ws.Range("D40").Value = "0.994"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").Value = "3.38"
$ws.Range("E41").Value = "  +6.87%  "
$ws.Range("D42").Value = "0.591"
$ws.Range("E42").Value = "  +5.33%  "
$ws.Range("D43").Value = "0.0538"
$ws.Range("E43").Value = "  +7.28%  "
$ws.Range("D44").Value = "1.25"
$ws.Range("E44").Value = "  +8.96%  "
$ws.Range("D45").Value = "10.16"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("D46").Value = "252.10"
$ws.Range("E46").Value = "  +30.62%  "
$ws.Range("E47").Value = "  +14.50%  "
$ws.Range("D48").Value = "0.0883"
$ws.Range("E48").Value = "  +9.20%  "
$ws.Range("E49").Value = "  +7.56%  "
$ws.Range("D50").Value = "1.885.09"
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("D51").Value = "16.75"
$ws.Range("E51").Value = "  +6.46%  "

Write-Host "Updated cryptos list"
